{"js": "// Replace the three-digit \u00f7 one-digit division problems in the table\n// with the new set of problems (text-for-text, 1:1 positional swap).\nconst replacements = [\n  [\"679\u00f75=\", \"459\u00f73=\"],\n  [\"562\u00f79=\", \"203\u00f77=\"],\n  [\"798\u00f76=\", \"637\u00f74=\"],\n  [\"388\u00f78=\", \"568\u00f73=\"],\n  [\"973\u00f74=\", \"596\u00f75=\"],\n  [\"702\u00f77=\", \"822\u00f79=\"],\n  [\"328\u00f74=\", \"533\u00f77=\"],\n  [\"195\u00f79=\", \"125\u00f78=\"],\n  [\"625\u00f75=\", \"372\u00f75=\"],\n  [\"650\u00f77=\", \"726\u00f73=\"],\n  [\"657\u00f79=\", \"600\u00f72=\"],\n  [\"674\u00f73=\", \"990\u00f74=\"],\n  [\"921\u00f76=\", \"861\u00f74=\"],\n  [\"640\u00f77=\", \"656\u00f72=\"],\n  [\"351\u00f73=\", \"950\u00f76=\"],\n  [\"903\u00f72=\", \"726\u00f79=\"],\n  [\"697\u00f79=\", \"179\u00f75=\"],\n  [\"396\u00f72=\", \"830\u00f75=\"],\n  [\"932\u00f72=\", \"344\u00f75=\"],\n  [\"285\u00f75=\", \"439\u00f78=\"],\n  [\"178\u00f78=\", \"146\u00f75=\"],\n  [\"372\u00f78=\", \"365\u00f76=\"],\n  [\"905\u00f75=\", \"182\u00f77=\"],\n  [\"187\u00f77=\", \"512\u00f77=\"],\n  [\"119\u00f76=\", \"723\u00f79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit \u00f7 one-digit division problems in the table\n# with the new set of problems (text-for-text, 1:1 positional swap).\n$pairs = @(\n    @(\"679\u00f75=\", \"459\u00f73=\"),\n    @(\"562\u00f79=\", \"203\u00f77=\"),\n    @(\"798\u00f76=\", \"637\u00f74=\"),\n    @(\"388\u00f78=\", \"568\u00f73=\"),\n    @(\"973\u00f74=\", \"596\u00f75=\"),\n    @(\"702\u00f77=\", \"822\u00f79=\"),\n    @(\"328\u00f74=\", \"533\u00f77=\"),\n    @(\"195\u00f79=\", \"125\u00f78=\"),\n    @(\"625\u00f75=\", \"372\u00f75=\"),\n    @(\"650\u00f77=\", \"726\u00f73=\"),\n    @(\"657\u00f79=\", \"600\u00f72=\"),\n    @(\"674\u00f73=\", \"990\u00f74=\"),\n    @(\"921\u00f76=\", \"861\u00f74=\"),\n    @(\"640\u00f77=\", \"656\u00f72=\"),\n    @(\"351\u00f73=\", \"950\u00f76=\"),\n    @(\"903\u00f72=\", \"726\u00f79=\"),\n    @(\"697\u00f79=\", \"179\u00f75=\"),\n    @(\"396\u00f72=\", \"830\u00f75=\"),\n    @(\"932\u00f72=\", \"344\u00f75=\"),\n    @(\"285\u00f75=\", \"439\u00f78=\"),\n    @(\"178\u00f78=\", \"146\u00f75=\"),\n    @(\"372\u00f78=\", \"365\u00f76=\"),\n    @(\"905\u00f75=\", \"182\u00f77=\"),\n    @(\"187\u00f77=\", \"512\u00f77=\"),\n    @(\"119\u00f76=\", \"723\u00f79=\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2)\n}\n"}
